$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2 through 420)
# from serial date 45175 to 45177, matching the upstream data refresh.
$ws.Range("C2:C420").Value = 45177
